$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 does not exist yet in the workbook. It needs the same per-column
# cell styles as the already-present (blank) template row 6, minus the
# trailing K column. Copy just the formatting (not values) from A6:J6
# down into A5:J5 so the new row matches the sheet's existing look.
$ws.Range("A6:J6").Copy()
$ws.Range("A5:J5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the new "Crumpet" facility row (row 5).
$ws.Range("A5").Value = "Crumpet GEF"
$ws.Range("B5").Value = 20001371
$ws.Range("C5").Value = "Crumpet exporter"
$ws.Range("D5").Value = "GBP"
$ws.Range("E5").Value = 7000000
$ws.Range("F5").Value = 3938753.8
$ws.Range("G5").Value = 777
$ws.Range("H5").Value = 456
$ws.Range("I5").Value = "GBP"
$ws.Range("J5").Value = "GBP"

# Fill in the new "Scone" facility row (row 6). Row 6 already existed as a
# blank, styled template row, so only the values need to be set; the
# trailing K6 cell remains blank exactly as before.
$ws.Range("A6").Value = "Scone GEF"
$ws.Range("B6").Value = 20001371
$ws.Range("C6").Value = "Scone exporter"
$ws.Range("D6").Value = "GBP"
$ws.Range("E6").Value = 770000
$ws.Range("F6").Value = 761579.37
$ws.Range("G6").Value = 777
$ws.Range("H6").Value = 456.77
$ws.Range("I6").Value = "GBP"
$ws.Range("J6").Value = "GBP"

# Reflect the user's final on-screen selection: cursor resting on the new
# data (A5:J6) rather than the old G2 selection scrolled to column B.
$ws.Range("A5:J6").Select()
